$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column, reusing the existing header formatting (same style as G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill data rows for the new column with 0
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
